$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5.830899999999999
$ws.Range("A12").Value = -21.412
$ws.Range("D14").Value = -8.2766
$ws.Range("D19").Value = -8.3384
$ws.Range("B23").Value = 8.845299999999996
$ws.Range("D24").Value = -7.861400000000003
$ws.Range("A27").Value = -21.92130000000001
$ws.Range("B28").Value = 5.890599999999999
$ws.Range("A32").Value = -21.009
$ws.Range("B32").Value = 6.0295
$ws.Range("B34").Value = 9.7142
$ws.Range("A36").Value = -19.9036
$ws.Range("A38").Value = -20.2983
$ws.Range("D38").Value = -7.340799999999999
$ws.Range("D41").Value = -8.370499999999995
$ws.Range("B42").Value = 10.2769
$ws.Range("A46").Value = -22.03660000000002
$ws.Range("B49").Value = 5.249200000000003
$ws.Range("D52").Value = -7.805500000000003
$ws.Range("A54").Value = -21.83570000000002
$ws.Range("B54").Value = 5.393499999999994
$ws.Range("A55").Value = -22.03180000000001
$ws.Range("A56").Value = -21.93860000000001
$ws.Range("A67").Value = -21.60869999999996
$ws.Range("A69").Value = -21.64569999999997
$ws.Range("A72").Value = -21.9051
$ws.Range("D72").Value = -7.544300000000004
$ws.Range("B78").Value = 8.803500000000003
$ws.Range("D78").Value = -7.968100000000003
$ws.Range("B80").Value = 9.701600000000003
$ws.Range("A83").Value = -21.61209999999999
$ws.Range("D83").Value = -7.842400000000002
$ws.Range("D85").Value = -8.811099999999993
$ws.Range("A86").Value = -21.7722
$ws.Range("D86").Value = -8.769000000000004
$ws.Range("D90").Value = -6.760099999999995
$ws.Range("A91").Value = -20.99509999999998
$ws.Range("A93").Value = -21.61340000000001
$ws.Range("D96").Value = -8.647999999999991
$ws.Range("B97").Value = 6.314099999999997
$ws.Range("A99").Value = -21.77979999999999
$ws.Range("B99").Value = 5.758699999999997
$ws.Range("B101").Value = 4.735699999999998
$ws.Range("D103").Value = -7.996399999999999
$ws.Range("A104").Value = -21.60169999999999
